# Convert the "HS" (Highest Score) column values that were stored as
# text like "113*" into plain numeric values (113), dropping the
# "not out" asterisk marker, for all the rows where that applied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "F2"  = 113
    "F3"  = 108
    "F8"  = 80
    "F9"  = 88
    "F11" = 89
    "F13" = 105
    "F15" = 71
    "F16" = 107
    "F17" = 58
    "F18" = 102
    "F23" = 108
    "F30" = 56
    "F31" = 32
    "F37" = 37
    "F38" = 27
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}

# Scroll the view so row 36 is at the top and select cell F2, matching
# the author's saved view state.
$ws.Range("F2").Select()
$excel.ActiveWindow.ScrollRow = 36
